$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (41) down into
# the two new rows (42 and 43) so the new cells inherit the same styles
# (bold/border index column, date-time number format column, etc.)
$ws.Range("A41:V41").Copy()
$ws.Range("A42:V43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 42: Dacia Buiucani vs Sparta Selemet
$ws.Cells.Item(42,1).Value = 41
$ws.Cells.Item(42,2).Value = "moldova"
$ws.Cells.Item(42,3).Value = "super-liga"
$ws.Cells.Item(42,4).Value = "2023-2024"
$ws.Cells.Item(42,5).Value = 45234.5
$ws.Cells.Item(42,6).Value = "Dacia Buiucani"
$ws.Cells.Item(42,7).Value = 1
$ws.Cells.Item(42,8).Value = "Sparta Selemet"
$ws.Cells.Item(42,9).Value = 1
$ws.Cells.Item(42,10).Value = 1.64
$ws.Cells.Item(42,11).Value = "03/11/2023 00:12"
$ws.Cells.Item(42,12).Value = 1.45
$ws.Cells.Item(42,13).Value = "04/11/2023 11:40"
$ws.Cells.Item(42,14).Value = 3.55
$ws.Cells.Item(42,15).Value = "03/11/2023 00:12"
$ws.Cells.Item(42,16).Value = 4.84
$ws.Cells.Item(42,17).Value = "04/11/2023 11:49"
$ws.Cells.Item(42,18).Value = 3.91
$ws.Cells.Item(42,19).Value = "03/11/2023 00:12"
$ws.Cells.Item(42,20).Value = 4.78
$ws.Cells.Item(42,21).Value = "04/11/2023 11:49"
$ws.Cells.Item(42,22).Value = "https://www.betexplorer.com/football/moldova/super-liga/dacia-buiucani-sparta-selemet/nZoCN2Y1/"

# Row 43: Petrocub vs Floresti
$ws.Cells.Item(43,1).Value = 42
$ws.Cells.Item(43,2).Value = "moldova"
$ws.Cells.Item(43,3).Value = "super-liga"
$ws.Cells.Item(43,4).Value = "2023-2024"
$ws.Cells.Item(43,5).Value = 45234.66666666666
$ws.Cells.Item(43,6).Value = "Petrocub"
$ws.Cells.Item(43,7).Value = 6
$ws.Cells.Item(43,8).Value = "Floresti"
$ws.Cells.Item(43,9).Value = 0
$ws.Cells.Item(43,10).Value = 1.33
$ws.Cells.Item(43,11).Value = "03/11/2023 04:13"
$ws.Cells.Item(43,12).Value = 1.29
$ws.Cells.Item(43,13).Value = "04/11/2023 15:52"
$ws.Cells.Item(43,14).Value = 4.25
$ws.Cells.Item(43,15).Value = "03/11/2023 04:13"
$ws.Cells.Item(43,16).Value = 5.28
$ws.Cells.Item(43,17).Value = "04/11/2023 15:56"
$ws.Cells.Item(43,18).Value = 6.11
$ws.Cells.Item(43,19).Value = "03/11/2023 04:13"
$ws.Cells.Item(43,20).Value = 7.25
$ws.Cells.Item(43,21).Value = "04/11/2023 15:56"
$ws.Cells.Item(43,22).Value = "https://www.betexplorer.com/football/moldova/super-liga/petrocub-hincesti-floresti/j7h3PO3k/"
